# Update "想去人数" (interested-count) figures in the 展览 and 全部类型 sheets
# F2: 326 -> 329
# F3: 1339 -> 1344
# F5: 67 -> 68

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 329
    $ws.Range("F3").Value = 1344
    $ws.Range("F5").Value = 68
}
